# Tutorial 6 attendance sheet update:
#  - Reformat the Date column from dd/mm/yyyy to dd-mm-yyyy (kept as literal text)
#  - Recompute the Total/Real/Invalid/Absent flags for several rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 3;  Date = "28-07-2022"; D = 1; E = 0; G = 1; H = 1 },
    @{ Row = 4;  Date = "01-08-2022"; D = 1; E = 1; G = 0; H = 0 },
    @{ Row = 5;  Date = "04-08-2022"; D = 1; E = 1; G = 0; H = 0 },
    @{ Row = 6;  Date = "08-08-2022"; D = 1; E = 1; G = 0; H = 0 },
    @{ Row = 7;  Date = "11-08-2022"; D = 0; E = 0; G = 0; H = 1 },
    @{ Row = 8;  Date = "15-08-2022"; D = 0; E = 0; G = 0; H = 1 },
    @{ Row = 9;  Date = "18-08-2022"; D = 0; E = 0; G = 0; H = 1 },
    @{ Row = 10; Date = "22-08-2022"; D = 1; E = 1; G = 0; H = 0 },
    @{ Row = 11; Date = "25-08-2022"; D = 1; E = 1; G = 0; H = 0 },
    @{ Row = 12; Date = "29-08-2022"; D = 1; E = 1; G = 0; H = 0 },
    @{ Row = 13; Date = "01-09-2022"; D = 0; E = 0; G = 0; H = 1 },
    @{ Row = 14; Date = "05-09-2022"; D = 0; E = 0; G = 0; H = 1 },
    @{ Row = 15; Date = "08-09-2022"; D = 0; E = 0; G = 0; H = 1 },
    @{ Row = 16; Date = "12-09-2022"; D = 0; E = 0; G = 0; H = 1 },
    @{ Row = 17; Date = "15-09-2022"; D = 0; E = 0; G = 0; H = 1 },
    @{ Row = 18; Date = "19-09-2022"; D = 0; E = 0; G = 0; H = 1 },
    @{ Row = 19; Date = "22-09-2022"; D = 0; E = 0; G = 0; H = 1 },
    @{ Row = 20; Date = "26-09-2022"; D = 0; E = 0; G = 0; H = 1 },
    @{ Row = 21; Date = "29-09-2022"; D = 0; E = 0; G = 0; H = 1 }
)

foreach ($r in $rows) {
    # Force the "Date" cell to stay plain text (otherwise Excel would
    # silently reinterpret the dd-mm-yyyy-looking string as a real date).
    $ws.Range("A" + $r.Row).NumberFormat = "@"
    $ws.Range("A" + $r.Row).Value = $r.Date

    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $r.E
    $ws.Range("G" + $r.Row).Value = $r.G
    $ws.Range("H" + $r.Row).Value = $r.H
}
